$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: date, header text, drop the old "CMH" cell, widen the merge ---
$ws.Range("A1").Value = 43776
$ws.Range("B1").Value = "Hashboard Settings"

# Stash B1's current look on a scratch cell, then unmerge/clear formats
# across B1:D1 so re-merging doesn't stamp B1's non-default style onto the
# newly-covered C1/D1 (which would otherwise leave stray blank <c> entries
# behind), and finally restore the look onto the merged anchor cell only.
$ws.Range("H1").Value = "scratch"
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

$ws.Range("B1:C1").UnMerge()
$ws.Range("B1:D1").ClearFormats()
$ws.Range("B1:D1").Merge()

$ws.Range("H1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("H1").Clear()

# --- Row 2: header labels (C2/D2 swap meaning: 14nm <-> Clarke) ---
$ws.Range("A2").Value = "Power Price ($/MW)"
$ws.Range("B2").Value = "16nm"
$ws.Range("C2").Value = "14nm"
$ws.Range("D2").Value = "Clarke"

# --- Row 3 ---
$ws.Range("A3").Value = "< 30"
$ws.Range("B3").Value = "Full"
$ws.Range("C3").Value = "Full"
$ws.Range("D3").Value = "Full"

# --- Row 4 ---
$ws.Range("A4").Value = "32 - 53"
$ws.Range("B4").Value = "Eco"
$ws.Range("C4").Value = "Full"
$ws.Range("D4").Value = "Full"

# --- Row 5 ---
$ws.Range("A5").Value = "53 - 71"
$ws.Range("B5").Value = "Eco"
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = "Full"

# --- Row 6 ---
$ws.Range("A6").Value = "71 - 92"
$ws.Range("B6").Value = "Standby"
$ws.Range("C6").Value = 29
$ws.Range("D6").Value = "Eco"

# --- Row 7 ---
$ws.Range("A7").Value = "> 92"
$ws.Range("B7").Value = "Standby"
$ws.Range("C7").Value = 29
$ws.Range("D7").Value = "Standby"
